$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on Price cells whose new values look numeric,
# so Excel stores them as text (matching the source data) instead of
# auto-converting them to actual numbers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.664.37"
$ws.Range("E2").Value = "  -2.11%  "

$ws.Range("D3").Value = "1.760.64"
$ws.Range("E3").Value = "  -3.03%  "

$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "324.15"
$ws.Range("E5").Value = "  -1.26%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "0.4308"
$ws.Range("E7").Value = "  -1.48%  "

$ws.Range("D8").Value = "0.3605"
$ws.Range("E8").Value = "  -1.73%  "

$ws.Range("D9").Value = "0.07566"
$ws.Range("E9").Value = "  -1.51%  "

$ws.Range("D10").Value = "42.20"
$ws.Range("E10").Value = "  -6.18%  "

$ws.Range("D11").Value = "1.110"
$ws.Range("E11").Value = "  -2.70%  "

$ws.Range("E12").Value = "  -0.01%  "

$ws.Range("D13").Value = "20.77"

$ws.Range("D14").Value = "6.073"
$ws.Range("E14").Value = "  -3.73%  "

$ws.Range("D15").Value = "7.242"
$ws.Range("E15").Value = "  -3.97%  "

$ws.Range("D16").Value = "1.758.91"
$ws.Range("E16").Value = "  -4.16%  "

$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("E18").Value = "  -1.39%  "

$ws.Range("D19").Value = "0.06429"
$ws.Range("E19").Value = "  -1.70%  "

$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").Value = "17.10"
$ws.Range("E21").Value = "  -2.36%  "

$ws.Range("D22").Value = "5.882"
$ws.Range("E22").Value = "  -6.07%  "

$ws.Range("D23").Value = "27.705.58"
$ws.Range("E23").Value = "  -2.08%  "

$ws.Range("E24").Value = "  -3.34%  "

$ws.Range("D25").Value = "2.086"
$ws.Range("E25").Value = "  +1.79%  "

$ws.Range("D26").Value = "162.52"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").Value = "20.59"
$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("D28").Value = "1.958.73"
$ws.Range("E28").Value = "  -3.81%  "

$ws.Range("D29").Value = "2.148"
$ws.Range("E29").Value = "  -6.58%  "

$ws.Range("D30").Value = "125.79"
$ws.Range("E30").Value = "  -2.31%  "

$ws.Range("D31").Value = "1.099"
$ws.Range("E31").Value = "  -9.66%  "

$ws.Range("D32").Value = "3.683"
$ws.Range("E32").Value = "  +5.39%  "

$ws.Range("D33").Value = "5.597"
$ws.Range("E33").Value = "  -6.16%  "

$ws.Range("D34").Value = "0.08964"
$ws.Range("E34").Value = "  -2.62%  "

$ws.Range("D35").Value = "12.21"
$ws.Range("E35").Value = "  -5.70%  "

$ws.Range("D36").Value = "0.02302"
$ws.Range("E36").Value = "  -1.96%  "

$ws.Range("D37").Value = "0.2117"
$ws.Range("E37").Value = "  -2.92%  "

$ws.Range("D38").Value = "0.06010"
$ws.Range("E38").Value = "  -3.09%  "

$ws.Range("D39").Value = "0.6358"
$ws.Range("E39").Value = "  -3.34%  "

$ws.Range("D40").Value = "4.959"
$ws.Range("E40").Value = "  -4.64%  "

$ws.Range("D41").Value = "1.191"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("E43").Value = "  -2.56%  "

$ws.Range("D44").Value = "7.897"
$ws.Range("E44").Value = "  -2.88%  "

$ws.Range("D45").Value = "13.33"
$ws.Range("E45").Value = "  -3.95%  "

$ws.Range("D46").Value = "0.5932"
$ws.Range("E46").Value = "  -3.10%  "

$ws.Range("D47").Value = "3.713"
$ws.Range("E47").Value = "  -1.13%  "

$ws.Range("D48").Value = "1.988"
$ws.Range("E48").Value = "  -1.73%  "

$ws.Range("D49").Value = "123.06"
$ws.Range("E49").Value = "  -2.19%  "

$ws.Range("D50").Value = "1.175"
$ws.Range("E50").Value = "  +1.44%  "

$ws.Range("D51").Value = "0.06876"
$ws.Range("E51").Value = "  -1.85%  "
